$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "97.215.86"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +4.83%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "3.129.86"
$c.Style = "Normal"
$ws.Range("E3").Value = "  +0.61%  "
$ws.Range("E4").Value = "  +0.04%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "241.12"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +2.40%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "611.88"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("E7").Value = "  +1.71%  "
$ws.Range("E8").Value = "  -1.05%  "
$ws.Range("E9").Value = "  +0.10%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "3.129.10"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.70%  "
$ws.Range("E11").Value = "  -0.26%  "
$ws.Range("E12").Value = "  -0.10%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "96.851.93"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +4.71%  "
$ws.Range("E14").Value = "  -1.36%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "33.94"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +0.17%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "5.45"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +0.58%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "3.712.33"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.79%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "3.127.41"
$c.Style = "Normal"
$ws.Range("E18").Value = "  +0.35%  "
$ws.Range("B19").Value = "BitcoinCash"
$ws.Range("C19").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "530.42"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +20.77%  "
$ws.Range("B20").Value = "SuiNetwork"
$ws.Range("C20").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "3.52"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -7.30%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "14.59"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.57%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "5.69"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -2.42%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.0000193"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -4.74%  "
$ws.Range("E24").Value = "  -2.69%  "
$ws.Range("E25").Value = "  +4.17%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "5.47"
$c.Style = "Normal"
$ws.Range("E26").Value = "  -1.95%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "11.63"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +0.75%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "3.297.48"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.70%  "
$ws.Range("E29").Value = "  -0.18%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "0.237"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +0.70%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.126"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("B32").Value = "Cronos"
$ws.Range("C32").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "0.175"
$c.Style = "Normal"
$ws.Range("E32").Value = "  -3.90%  "
$ws.Range("B33").Value = "Binance-PegBSC-USD"
$ws.Range("C33").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.67%  "
$ws.Range("B34").Value = "InternetComputer(DFINITY)"
$ws.Range("C34").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "9.01"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -1.59%  "
$ws.Range("B35").Value = "EthereumClassic"
$ws.Range("C35").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "26.75"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +4.16%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.152"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -7.40%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "7.33"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -9.13%  "
$ws.Range("B38").Value = "PancakeSwap"
$ws.Range("C38").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.89"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.61%  "
$ws.Range("B39").Value = "Bittensor"
$ws.Range("C39").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "486.26"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +4.82%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "24.26"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.44%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.442"
$c.Style = "Normal"
$ws.Range("E41").Value = "  +3.40%  "
$ws.Range("E42").Value = "  -4.39%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "3.59"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -10.32%  "
$ws.Range("E44").Value = "  -0.01%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "3.18"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -4.83%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "162.14"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.62%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "0.702"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +2.95%  "
$ws.Range("E48").Value = "  +4.41%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "4.49"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +3.37%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "44.36"
$c.Style = "Normal"
$ws.Range("E50").Value = "  +1.27%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.999"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +0.07%  "
